$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (37 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6802533.5
$ws.Range("J17").Value = 7288293
$ws.Range("L17").Value = 21864879
$ws.Range("N17").Value = -21865215
$ws.Range("H53").Value = 3792.6365
$ws.Range("J53").Value = 5878.143
$ws.Range("L53").Value = 5878.143
$ws.Range("N53").Value = -7152.143
$ws.Range("H116").Value = 5307.9165
$ws.Range("J116").Value = 6350.75
$ws.Range("L116").Value = 6350.75
$ws.Range("N116").Value = -13234.75
$ws.Range("H129").Value = 189686.1
$ws.Range("J129").Value = 197107.12
$ws.Range("L129").Value = 591321.36
$ws.Range("N129").Value = -601321.36
$ws.Range("H132").Value = 2632.6
$ws.Range("I132").Value = 3060.0967
$ws.Range("J132").Value = 1160.1111
$ws.Range("K132").Value = 9180.2901
$ws.Range("L132").Value = 3480.3333
$ws.Range("M132").Value = -6650.2901
$ws.Range("N132").Value = -8540.3333
$ws.Range("H135").Value = 17242502
$ws.Range("I135").Value = 848.13043
$ws.Range("J135").Value = 83335510
$ws.Range("K135").Value = 7633.173870000001
$ws.Range("L135").Value = 750019590
$ws.Range("M135").Value = -5098.173870000001
$ws.Range("N135").Value = -750024660
$ws.Range("H137").Value = 1231.5366
$ws.Range("I137").Value = 1238.5555
$ws.Range("J137").Value = 1218
$ws.Range("K137").Value = 3715.6665
$ws.Range("L137").Value = 3654
$ws.Range("M137").Value = -1165.6665
$ws.Range("N137").Value = -8754

# --- Sheet: ARM (47 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6226.4165
$ws.Range("I32").Value = 4923.8833
$ws.Range("J32").Value = 12739.083
$ws.Range("K32").Value = 4923.8833
$ws.Range("L32").Value = 12739.083
$ws.Range("M32").Value = -4636.8833
$ws.Range("N32").Value = -13313.083
$ws.Range("H61").Value = 4455.263
$ws.Range("I61").Value = 4617.857
$ws.Range("K61").Value = 4617.857
$ws.Range("M61").Value = -4405.857
$ws.Range("H74").Value = 26317244
$ws.Range("I74").Value = 40000620
$ws.Range("J74").Value = 3061.4614
$ws.Range("K74").Value = 40000620
$ws.Range("L74").Value = 3061.4614
$ws.Range("M74").Value = -39999746
$ws.Range("N74").Value = -4809.4614
$ws.Range("H77").Value = 26317244
$ws.Range("I77").Value = 40000620
$ws.Range("J77").Value = 3061.4614
$ws.Range("K77").Value = 200003100
$ws.Range("L77").Value = 15307.307
$ws.Range("M77").Value = -199998732
$ws.Range("N77").Value = -24043.307
$ws.Range("H102").Value = 1539.3334
$ws.Range("I102").Value = 1397.6666
$ws.Range("J102").Value = 1751.8334
$ws.Range("K102").Value = 1397.6666
$ws.Range("L102").Value = 1751.8334
$ws.Range("M102").Value = 224.3334
$ws.Range("N102").Value = -4995.8334
$ws.Range("H132").Value = 13718.568
$ws.Range("I132").Value = 2173.6765
$ws.Range("J132").Value = 52971.2
$ws.Range("K132").Value = 6521.029500000001
$ws.Range("L132").Value = 158913.6
$ws.Range("M132").Value = -3991.029500000001
$ws.Range("N132").Value = -163973.6
$ws.Range("H136").Value = 4455.263
$ws.Range("I136").Value = 4617.857
$ws.Range("K136").Value = 13853.571
$ws.Range("M136").Value = -11303.571
$ws.Range("H138").Value = 50179
$ws.Range("J138").Value = 50179
$ws.Range("L138").Value = 50179
$ws.Range("N138").Value = -60459

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1112.4286
$ws.Range("I99").Value = 1207.5
$ws.Range("J99").Value = 874.75
$ws.Range("K99").Value = 1207.5
$ws.Range("L99").Value = 874.75
$ws.Range("M99").Value = 290.5
$ws.Range("N99").Value = -3870.75
$ws.Range("H134").Value = 3375.1428
$ws.Range("I134").Value = 3444.5144
$ws.Range("K134").Value = 10333.5432
$ws.Range("M134").Value = -7798.5432

# --- Sheet: CRP (47 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 351.4
$ws.Range("I7").Value = 74
$ws.Range("J7").Value = 628.8
$ws.Range("K7").Value = 74
$ws.Range("L7").Value = 628.8
$ws.Range("M7").Value = 39
$ws.Range("N7").Value = -854.8
$ws.Range("H22").Value = 175.8
$ws.Range("I22").Value = 175.8
$ws.Range("K22").Value = 175.8
$ws.Range("M22").Value = 174.2
$ws.Range("H31").Value = 3653.8948
$ws.Range("I31").Value = 2882
$ws.Range("J31").Value = 4511.5557
$ws.Range("K31").Value = 2882
$ws.Range("L31").Value = 4511.5557
$ws.Range("M31").Value = -2587
$ws.Range("N31").Value = -5101.5557
$ws.Range("H34").Value = 3653.8948
$ws.Range("I34").Value = 2882
$ws.Range("J34").Value = 4511.5557
$ws.Range("K34").Value = 2882
$ws.Range("L34").Value = 4511.5557
$ws.Range("M34").Value = -2680
$ws.Range("N34").Value = -4915.5557
$ws.Range("H99").Value = 17860184
$ws.Range("I99").Value = 2620.65
$ws.Range("J99").Value = 62504090
$ws.Range("K99").Value = 2620.65
$ws.Range("L99").Value = 62504090
$ws.Range("M99").Value = -1122.65
$ws.Range("N99").Value = -62507086
$ws.Range("H126").Value = 17860184
$ws.Range("I126").Value = 2620.65
$ws.Range("J126").Value = 62504090
$ws.Range("K126").Value = 7861.950000000001
$ws.Range("L126").Value = 187512270
$ws.Range("M126").Value = -5391.950000000001
$ws.Range("N126").Value = -187517210
$ws.Range("H132").Value = 3761.8333
$ws.Range("I132").Value = 2669.3076
$ws.Range("K132").Value = 8007.9228
$ws.Range("M132").Value = -5477.9228
$ws.Range("H134").Value = 1340
$ws.Range("I134").Value = 1124.6154
$ws.Range("K134").Value = 3373.8462
$ws.Range("M134").Value = -838.8462

# --- Sheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1535.0555
$ws.Range("I5").Value = 687.7
$ws.Range("J5").Value = 2594.25
$ws.Range("K5").Value = 2063.1
$ws.Range("L5").Value = 7782.75
$ws.Range("M5").Value = -1951.1
$ws.Range("N5").Value = -8006.75
$ws.Range("H135").Value = 1535.0555
$ws.Range("I135").Value = 687.7
$ws.Range("J135").Value = 2594.25
$ws.Range("K135").Value = 6189.3
$ws.Range("L135").Value = 23348.25
$ws.Range("M135").Value = -3654.3
$ws.Range("N135").Value = -28418.25

# --- Sheet: LTW (12 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1310901.5
$ws.Range("J122").Value = 4159.875
$ws.Range("L122").Value = 12479.625
$ws.Range("N122").Value = -17379.625
$ws.Range("H132").Value = 417336.78
$ws.Range("I132").Value = 603699.5
$ws.Range("K132").Value = 1811098.5
$ws.Range("M132").Value = -1808568.5
$ws.Range("H136").Value = 1283.8438
$ws.Range("I136").Value = 1192.3334
$ws.Range("K136").Value = 3577.0002
$ws.Range("M136").Value = -1027.0002

# --- Sheet: WVR (4 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 351.36365
$ws.Range("I100").Value = 357.22223
$ws.Range("K100").Value = 714.44446
$ws.Range("M100").Value = -173.44446

Write-Host "Applied all Typhon_Profits.xlsx market-data updates"